$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 129, pushing the existing rows 129:150 down to 130:151
# (this matches the diff, which shows every row from 129..150 taking on the
# values previously held by the row above it, and a brand-new row 151
# appearing with the data that used to live in row 150).
$ws.Rows(129).Insert()

# Populate the newly inserted row 129 with the new weekly entry.
$ws.Cells.Item(129, 1).Value = 10
$ws.Cells.Item(129, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(129, 3).Value = "La Araucanía"
$ws.Cells.Item(129, 4).Value = 44505
$ws.Cells.Item(129, 5).Value = 9
$ws.Cells.Item(129, 6).Value = 100112005
$ws.Cells.Item(129, 7).Value = "Puerro"
$ws.Cells.Item(129, 8).Value = "Azul de Maquehue"
$ws.Cells.Item(129, 9).Value = "Primera"
$ws.Cells.Item(129, 10).Value = 65
$ws.Cells.Item(129, 11).Value = 7000
$ws.Cells.Item(129, 12).Value = 7000
$ws.Cells.Item(129, 13).Value = 7000
$ws.Cells.Item(129, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(129, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(129, 16).Value = 583
$ws.Cells.Item(129, 17).Value = 12
$ws.Cells.Item(129, 18).Value = "Hortaliza"
